$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Swap the two zero-width DDE-link bookmarks' names (id 0 <-> id 1 keep
#    their ids/positions, but the names attached to them are exchanged).
# ---------------------------------------------------------------------------
$bmA = $d.Bookmarks("__DdeLink__113_1728794273")
$startA = $bmA.Start
$endA = $bmA.End

$bmB = $d.Bookmarks("__DdeLink__115_1728794273")
$startB = $bmB.Start
$endB = $bmB.End

$bmA.Delete()
$bmB.Delete()

# Re-create with swapped names. Because new bookmarks are inserted just
# before any existing bookmark that starts at the very same location, we
# add the "113" one first so the "115" one (added second) ends up first in
# the document -- matching id=0 -> 115 / id=1 -> 113 from the target.
$rngA = $d.Range($startA, $endA)
$rngA.Bookmarks.Add("__DdeLink__113_1728794273")
$rngB = $d.Range($startA, $endA)
$rngB.Bookmarks.Add("__DdeLink__115_1728794273")

# ---------------------------------------------------------------------------
# 2) Merge the split e-mail runs "Viievskiy.Anton@yandex." + "ua" back into
#    a single run with the full address (Find can locate, but can't Replace
#    across this particular hyperlink boundary, so splice the Range.Text).
# ---------------------------------------------------------------------------
$mailRng = $d.Content
$mailRng.Find.Execute("yandex") | Out-Null
$mailRng.Text = "yandex"

# ---------------------------------------------------------------------------
# 3) Insert " Foundation 6," into the skills sentence, splitting the run
#    into three runs that share identical formatting.
# ---------------------------------------------------------------------------
$skillsRng = $d.Content
$skillsRng.Find.Execute("Skeleton, ") | Out-Null
$skillsRng.Collapse(0)
$insertStart = $skillsRng.Start
$skillsRng.InsertAfter("Foundation 6, ")

# Temporarily bold just "Foundation 6," (not the trailing space) so the
# engine keeps it as its own run, then remove the bold again -- the run
# boundaries stay in place once formatting diverged and reconverged.
$newRng = $d.Range($insertStart, $insertStart + 13)
$newRng.Font.Bold = 1
$newRng.Font.Bold = 0

# ---------------------------------------------------------------------------
# 4) Flip the "Normal" style's overflowPunct paragraph setting from false
#    to true (VBA/COM: ParagraphFormat.HangingPunctuation).
# ---------------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = $true
